# Update cryptocurrency price/volume data per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (text). Column D holds numeric-looking strings
# (e.g. "1.001", "30.805.79") that must stay plain text, matching the source
# data which uses "." as both a thousands separator and a decimal point.
$updates = [ordered]@{
    'D2' = '30.805.79'
    'E2' = '  +0.64%  '
    'D3' = '1.888.23'
    'E3' = '  +1.10%  '
    'D4' = '1.001'
    'E4' = '  +0.01%  '
    'D5' = '239.83'
    'E5' = '  +2.00%  '
    'E6' = '  -0.02%  '
    'D7' = '0.4792'
    'D8' = '0.2951'
    'E8' = '  +6.65%  '
    'D9' = '0.06616'
    'E9' = '  +3.75%  '
    'D10' = '18.80'
    'E10' = '  +5.03%  '
    'D11' = '101.07'
    'E11' = '  +18.65%  '
    'D12' = '1.896.64'
    'E12' = '  +1.76%  '
    'D13' = '0.07637'
    'E13' = '  +2.36%  '
    'D14' = '5.133'
    'E14' = '  +3.00%  '
    'D15' = '0.6589'
    'E15' = '  +3.59%  '
    'D16' = '303.79'
    'E16' = '  +25.99%  '
    'D17' = '30.773.64'
    'E17' = '  +0.63%  '
    'D18' = '13.18'
    'E18' = '  +2.47%  '
    'E19' = '  +0.06%  '
    'D20' = '0.000007620'
    'E20' = '  +3.28%  '
    'D21' = '2.140.66'
    'E21' = '  +2.21%  '
    'D22' = '1.001'
    'E22' = '  -0.01%  '
    'D23' = '5.170'
    'E23' = '  +3.47%  '
    'D24' = '6.188'
    'E24' = '  +2.48%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D25' = '169.00'
    'E25' = '  +1.90%  '
    'B26' = 'Cosmos'
    'C26' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D26' = '9.325'
    'E26' = '  -0.76%  '
    'D27' = '20.70'
    'E27' = '  +13.67%  '
    'D28' = '1.965'
    'E28' = '  +3.69%  '
    'D29' = '0.1121'
    'E29' = '  +9.61%  '
    'D30' = '1.346'
    'E30' = '  -2.45%  '
    'D31' = '4.192'
    'E31' = '  +2.13%  '
    'D32' = '3.998'
    'E32' = '  +3.35%  '
    'D33' = '0.05083'
    'E33' = '  +2.93%  '
    'D34' = '0.7443'
    'E34' = '  +4.93%  '
    'D35' = '1.160'
    'E35' = '  +0.60%  '
    'D36' = '2.721'
    'E36' = '  +0.58%  '
    'D37' = '0.01987'
    'E37' = '  +4.05%  '
    'D38' = '2.709'
    'E38' = '  +0.75%  '
    'D39' = '2.053'
    'E39' = '  +2.71%  '
    'D40' = '109.50'
    'E40' = '  +3.46%  '
    'D41' = '0.8808'
    'E41' = '  -0.07%  '
    'E42' = '  +0.02%  '
    'D43' = '0.4203'
    'E43' = '  +2.19%  '
    'D44' = '5.655'
    'E44' = '  +1.81%  '
    'D45' = '67.94'
    'E45' = '  +8.23%  '
    'D46' = '7.366'
    'E46' = '  -0.67%  '
    'D47' = '9.128'
    'E47' = '  +5.18%  '
    'D48' = '0.1230'
    'E48' = '  -0.07%  '
    'D49' = '34.90'
    'E49' = '  +3.43%  '
    'D50' = '0.05656'
    'E50' = '  +1.50%  '
    'D51' = '1.398'
    'E51' = '  +1.16%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    if ($cellRef[0] -eq "D") {
        # Force text storage so Excel does not coerce these into numbers
        # (which would strip trailing zeros / misparse the dotted values),
        # then drop back to the workbook default style (no NumberFormat
        # override survives), matching the original unstyled cells.
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$cellRef]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$cellRef]
    }
}
